# Changed date of sample flight
#
# The sample response row on "Form1" recorded a flight started on
# 2019-01-01 15:37:45 (serial 43466.65121527778). Move that sample date
# forward 9 days to 2019-01-10, keeping the same time-of-day, and leave
# the selection on the edited cell (B2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form1")
$ws.Activate()

$ws.Range("B2").Value = 43475.65121527778

$ws.Range("B2").Select()
